$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Forecast Comparison")

# --- 1. Insert a new column before "ASIN" (currently column B) for the new
#        "Week_Start_Date" field. This shifts ASIN .. is_holiday_week one
#        column to the right (B:I -> C:J).
$ws.Columns.Item(2).Insert()

# --- 2. New column header.
$ws.Range("B1").Value = "Week_Start_Date"

# --- 3. Populate the new column with the week-start date for each row.
#        These are plain text values (e.g. "2025-01-05"), not real Excel
#        dates, so force a text number format before assigning the value.
$weekStartDates = @(
    "2025-01-05",
    "2025-01-12",
    "2025-01-19",
    "2025-01-26",
    "2025-02-02",
    "2025-02-09",
    "2025-02-16",
    "2025-02-23",
    "2025-03-02",
    "2025-03-09",
    "2025-03-16",
    "2025-03-23",
    "2025-03-30",
    "2025-04-06",
    "2025-04-13",
    "2025-04-20"
)

for ($i = 0; $i -lt $weekStartDates.Length; $i++) {
    $row = $i + 2
    $cell = $ws.Cells.Item($row, 2)
    $cell.NumberFormat = "@"
    $cell.Value = $weekStartDates[$i]
    # Drop the temporary text-number-format now that the literal text is
    # locked in, so the cell is left with the default (no explicit) style.
    $cell.ClearFormats()
}

# --- 4. Correct the week labels in column A: drop the leading zero on the
#        single-digit weeks (W01 -> W1 ... W09 -> W9). W10-W16 already have
#        no leading zero and stay as-is.
$weekLabels = @{
    2  = "W1"
    3  = "W2"
    4  = "W3"
    5  = "W4"
    6  = "W5"
    7  = "W6"
    8  = "W7"
    9  = "W8"
    10 = "W9"
}

foreach ($row in $weekLabels.Keys) {
    $ws.Cells.Item($row, 1).Value = $weekLabels[$row]
}
